$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 59.2
$ws.Range("N2").Value = 85.8724807945396

$ws.Range("K3").Value = 56.2
$ws.Range("N3").Value = 85.8724807945396

$ws.Range("K4").Value = 54.8
$ws.Range("N4").Value = 85.8724807945396

$ws.Range("K5").Value = 51.8
$ws.Range("N5").Value = 85.8724807945396
